# Updated soil type data
# Adds a new "pro_usda_soil_order" column to the "profile" sheet, inserted
# immediately before the existing "pro_soil_taxon" column (column N / 14),
# and fills in the USDA soil order for the two existing profile rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("profile")

# Insert a new blank column at N (14); everything from N onward shifts right.
$ws.Columns.Item(14).Insert() | Out-Null

# Header + new data for the inserted column.
$ws.Range("N1").Value = "pro_usda_soil_order"
$ws.Range("N4").Value = "Inceptisols"
$ws.Range("N5").Value = "Spodosols"

# Reflect the author's final selection: cursor left on the new column in
# "profile", with "metadata" as the active/selected tab.
$ws.Activate() | Out-Null
$ws.Range("N6").Select() | Out-Null

$wsMeta = $wb.Worksheets.Item("metadata")
$wsMeta.Activate() | Out-Null
$wsMeta.Range("A4").Select() | Out-Null
